# Natmi following Dr Hou advice
# Updates the recomputed NATMI ligand-receptor edge statistics for rows 2-10
# (sheet1 / "Hbegf-Cd82") to reflect the refreshed ligand/receptor expressing
# cell counts (1 -> 3) and the downstream recomputed expression / specificity
# / edge-weight values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.739149
$ws.Range("H2").Value = 41.217447
$ws.Range("I2").Value = 0.6130043224686931
$ws.Range("J2").Value = 0.6130043224686931
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 27.493006
$ws.Range("N2").Value = 82.479018
$ws.Range("O2").Value = 0.6421229928833972
$ws.Range("P2").Value = 0.6421229928833972
$ws.Range("Q2").Value = 377.7305058918939
$ws.Range("R2").Value = 3399.574553027046
$ws.Range("S2").Value = 0.3936241701940564
$ws.Range("T2").Value = 0.3936241701940564
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.739149
$ws.Range("H3").Value = 41.217447
$ws.Range("I3").Value = 0.6130043224686931
$ws.Range("J3").Value = 0.6130043224686931
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.708541
$ws.Range("N3").Value = 5.125623
$ws.Range("O3").Value = 0.03990445644190354
$ws.Range("P3").Value = 0.03990445644190353
$ws.Range("Q3").Value = 23.473899371609
$ws.Range("R3").Value = 211.265094344481
$ws.Range("S3").Value = 0.02446160428465055
$ws.Range("T3").Value = 0.02446160428465055
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.739149
$ws.Range("H4").Value = 41.217447
$ws.Range("I4").Value = 0.6130043224686931
$ws.Range("J4").Value = 0.6130043224686931
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.61424733333333
$ws.Range("N4").Value = 40.842742
$ws.Range("O4").Value = 0.3179725506746993
$ws.Range("P4").Value = 0.3179725506746992
$ws.Range("Q4").Value = 187.0481726355193
$ws.Range("R4").Value = 1683.433553719674
$ws.Range("S4").Value = 0.1949185479899862
$ws.Range("T4").Value = 0.1949185479899862
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.034036666666666
$ws.Range("H5").Value = 12.10211
$ws.Range("I5").Value = 0.1799879973398545
$ws.Range("J5").Value = 0.1799879973398545
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 27.493006
$ws.Range("N5").Value = 82.479018
$ws.Range("O5").Value = 0.6421229928833972
$ws.Range("P5").Value = 0.6421229928833972
$ws.Range("Q5").Value = 110.9077942808866
$ws.Range("R5").Value = 998.1701485279799
$ws.Range("S5").Value = 0.1155744315349563
$ws.Range("T5").Value = 0.1155744315349563
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.034036666666666
$ws.Range("H6").Value = 12.10211
$ws.Range("I6").Value = 0.1799879973398545
$ws.Range("J6").Value = 0.1799879973398545
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.708541
$ws.Range("N6").Value = 5.125623
$ws.Range("O6").Value = 0.03990445644190354
$ws.Range("P6").Value = 0.03990445644190353
$ws.Range("Q6").Value = 6.892317040503333
$ws.Range("R6").Value = 62.03085336453
$ws.Range("S6").Value = 0.007182323199913675
$ws.Range("T6").Value = 0.007182323199913673
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.034036666666666
$ws.Range("H7").Value = 12.10211
$ws.Range("I7").Value = 0.1799879973398545
$ws.Range("J7").Value = 0.1799879973398545
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.61424733333333
$ws.Range("N7").Value = 40.842742
$ws.Range("O7").Value = 0.3179725506746993
$ws.Range("P7").Value = 0.3179725506746992
$ws.Range("Q7").Value = 54.92037293173555
$ws.Range("R7").Value = 494.28335638562
$ws.Range("S7").Value = 0.05723124260498453
$ws.Range("T7").Value = 0.05723124260498452
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.639623666666666
$ws.Range("H8").Value = 13.918871
$ws.Range("I8").Value = 0.2070076801914524
$ws.Range("J8").Value = 0.2070076801914524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 27.493006
$ws.Range("N8").Value = 82.479018
$ws.Range("O8").Value = 0.6421229928833972
$ws.Range("P8").Value = 0.6421229928833972
$ws.Range("Q8").Value = 127.5572013054086
$ws.Range("R8").Value = 1148.014811748678
$ws.Range("S8").Value = 0.1329243911543846
$ws.Range("T8").Value = 0.1329243911543846
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.639623666666666
$ws.Range("H9").Value = 13.918871
$ws.Range("I9").Value = 0.2070076801914524
$ws.Range("J9").Value = 0.2070076801914524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.708541
$ws.Range("N9").Value = 5.125623
$ws.Range("O9").Value = 0.03990445644190354
$ws.Range("P9").Value = 0.03990445644190353
$ws.Range("Q9").Value = 7.926987259070333
$ws.Range("R9").Value = 71.342885331633
$ws.Range("S9").Value = 0.008260528957339312
$ws.Range("T9").Value = 0.00826052895733931
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.639623666666666
$ws.Range("H10").Value = 13.918871
$ws.Range("I10").Value = 0.2070076801914524
$ws.Range("J10").Value = 0.2070076801914524
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.61424733333333
$ws.Range("N10").Value = 40.842742
$ws.Range("O10").Value = 0.3179725506746993
$ws.Range("P10").Value = 0.3179725506746992
$ws.Range("Q10").Value = 63.16498413158688
$ws.Range("R10").Value = 568.484857184282
$ws.Range("S10").Value = 0.06582276007972854
$ws.Range("T10").Value = 0.06582276007972852
